$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 107.75
$ws.Range("I12").Value = 135
$ws.Range("J12").Value = 80.5
$ws.Range("K12").Value = 135
$ws.Range("L12").Value = 80.5
$ws.Range("M12").Value = 35
$ws.Range("N12").Value = -420.5

$ws.Range("H86").Value = 1817.3572
$ws.Range("I86").Value = 992.875
$ws.Range("J86").Value = 2916.6667
$ws.Range("K86").Value = 992.875
$ws.Range("L86").Value = 2916.6667
$ws.Range("M86").Value = 130.125
$ws.Range("N86").Value = -5162.6667

$ws.Range("H89").Value = 1817.3572
$ws.Range("I89").Value = 992.875
$ws.Range("J89").Value = 2916.6667
$ws.Range("K89").Value = 4964.375
$ws.Range("L89").Value = 14583.3335
$ws.Range("M89").Value = 651.625
$ws.Range("N89").Value = -25815.3335

$ws.Range("H125").Value = 1731.1666
$ws.Range("I125").Value = 1705.5714
$ws.Range("J125").Value = 1767
$ws.Range("K125").Value = 15350.1426
$ws.Range("L125").Value = 15903
$ws.Range("M125").Value = -12890.1426
$ws.Range("N125").Value = -20823

$ws.Range("H132").Value = 23820958
$ws.Range("I132").Value = 37052830
$ws.Range("K132").Value = 111158490
$ws.Range("M132").Value = -111155960

$ws.Range("H137").Value = 1413.1923
$ws.Range("I137").Value = 950.4286
$ws.Range("J137").Value = 1953.0834
$ws.Range("K137").Value = 2851.2858
$ws.Range("L137").Value = 5859.2502
$ws.Range("M137").Value = -301.2857999999997
$ws.Range("N137").Value = -10959.2502

$ws.Range("H139").Value = 31080
$ws.Range("J139").Value = 31080
$ws.Range("L139").Value = 31080
$ws.Range("N139").Value = -41360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4093.39
$ws.Range("I32").Value = 3826.2246
$ws.Range("K32").Value = 3826.2246
$ws.Range("M32").Value = -3539.2246

$ws.Range("H61").Value = 2081
$ws.Range("I61").Value = 1774.6666
$ws.Range("K61").Value = 1774.6666
$ws.Range("M61").Value = -1562.6666

$ws.Range("H110").Value = 2269.4285
$ws.Range("I110").Value = 1869.875
$ws.Range("J110").Value = 2802.1667
$ws.Range("K110").Value = 1869.875
$ws.Range("L110").Value = 2802.1667
$ws.Range("M110").Value = 175.125
$ws.Range("N110").Value = -6892.1667

$ws.Range("H133").Value = 38520
$ws.Range("J133").Value = 38520
$ws.Range("L133").Value = 38520
$ws.Range("N133").Value = -43580

$ws.Range("H136").Value = 2081
$ws.Range("I136").Value = 1774.6666
$ws.Range("K136").Value = 5323.9998
$ws.Range("M136").Value = -2773.9998

$ws.Range("H139").Value = 38519.168
$ws.Range("J139").Value = 38519.168
$ws.Range("L139").Value = 38519.168
$ws.Range("N139").Value = -48799.168

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2363.25
$ws.Range("I20").Value = 1883.3334
$ws.Range("J20").Value = 3803
$ws.Range("K20").Value = 1883.3334
$ws.Range("L20").Value = 3803
$ws.Range("M20").Value = -1636.3334
$ws.Range("N20").Value = -4297

$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H105").Value = 66668984
$ws.Range("I105").Value = 76925140
$ws.Range("J105").Value = 3999.5
$ws.Range("K105").Value = 76925140
$ws.Range("L105").Value = 3999.5
$ws.Range("M105").Value = -76923393
$ws.Range("N105").Value = -7493.5

$ws.Range("H107").Value = 1241.7916
$ws.Range("I107").Value = 1054.2778
$ws.Range("J107").Value = 1804.3334
$ws.Range("K107").Value = 1054.2778
$ws.Range("L107").Value = 1804.3334
$ws.Range("M107").Value = 865.7221999999999
$ws.Range("N107").Value = -5644.3334

$ws.Range("H132").Value = 49714.145
$ws.Range("J132").Value = 49714.145
$ws.Range("L132").Value = 49714.145
$ws.Range("N132").Value = -59834.145

$ws.Range("H134").Value = 13078.637
$ws.Range("I134").Value = 1963
$ws.Range("J134").Value = 19430.428
$ws.Range("K134").Value = 5889
$ws.Range("L134").Value = 58291.284
$ws.Range("M134").Value = -3354
$ws.Range("N134").Value = -63361.284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 828.5
$ws.Range("I10").Value = 828.5
$ws.Range("K10").Value = 828.5
$ws.Range("M10").Value = -689.5

$ws.Range("H31").Value = 1649.3043
$ws.Range("I31").Value = 1273.2727
$ws.Range("K31").Value = 1273.2727
$ws.Range("M31").Value = -978.2727

$ws.Range("H34").Value = 1649.3043
$ws.Range("I34").Value = 1273.2727
$ws.Range("K34").Value = 1273.2727
$ws.Range("M34").Value = -1071.2727

$ws.Range("H105").Value = 734.9474
$ws.Range("I105").Value = 731.1539
$ws.Range("J105").Value = 743.1667
$ws.Range("K105").Value = 731.1539
$ws.Range("L105").Value = 743.1667
$ws.Range("M105").Value = 1015.8461
$ws.Range("N105").Value = -4237.1667

$ws.Range("H107").Value = 493.07693
$ws.Range("I107").Value = 443.75
$ws.Range("K107").Value = 443.75
$ws.Range("M107").Value = 1476.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1275.258
$ws.Range("I5").Value = 1378.1923
$ws.Range("J5").Value = 740
$ws.Range("K5").Value = 4134.5769
$ws.Range("L5").Value = 2220
$ws.Range("M5").Value = -4022.5769
$ws.Range("N5").Value = -2444

$ws.Range("H130").Value = 1777.5

$ws.Range("H131").Value = 1769.36
$ws.Range("J131").Value = 1798.0204
$ws.Range("L131").Value = 5394.0612
$ws.Range("N131").Value = -15474.0612

$ws.Range("H132").Value = 1666.6666
$ws.Range("I132").Value = 1250
$ws.Range("K132").Value = 11250
$ws.Range("M132").Value = -8720

$ws.Range("H135").Value = 1275.258
$ws.Range("I135").Value = 1378.1923
$ws.Range("J135").Value = 740
$ws.Range("K135").Value = 12403.7307
$ws.Range("L135").Value = 6660
$ws.Range("M135").Value = -9868.7307
$ws.Range("N135").Value = -11730

$ws.Range("H141").Value = 3666.6667
$ws.Range("I141").Value = 3250
$ws.Range("K141").Value = 9750
$ws.Range("M141").Value = -4570

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 679.4583
$ws.Range("I107").Value = 785.9286
$ws.Range("J107").Value = 530.4
$ws.Range("K107").Value = 785.9286
$ws.Range("L107").Value = 530.4
$ws.Range("M107").Value = 1134.0714
$ws.Range("N107").Value = -4370.4

$ws.Range("H113").Value = 1572.9166
$ws.Range("I113").Value = 1485.7142
$ws.Range("J113").Value = 1695
$ws.Range("K113").Value = 1485.7142
$ws.Range("L113").Value = 1695
$ws.Range("M113").Value = 684.2858000000001
$ws.Range("N113").Value = -6035

$ws.Range("H122").Value = 2021.2222
$ws.Range("I122").Value = 1845.4667
$ws.Range("J122").Value = 2900
$ws.Range("K122").Value = 5536.4001
$ws.Range("L122").Value = 8700
$ws.Range("M122").Value = -3086.4001
$ws.Range("N122").Value = -13600

$ws.Range("H132").Value = 3215
$ws.Range("I132").Value = 3416.1667
$ws.Range("K132").Value = 10248.5001
$ws.Range("M132").Value = -7718.500100000001

$ws.Range("H133").Value = 44599
$ws.Range("J133").Value = 44599
$ws.Range("L133").Value = 44599
$ws.Range("N133").Value = -54719

$ws.Range("H139").Value = 35000
$ws.Range("J139").Value = 35000
$ws.Range("L139").Value = 35000
$ws.Range("N139").Value = -45280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 24700.408
$ws.Range("I132").Value = 1514.3182
$ws.Range("K132").Value = 4542.9546
$ws.Range("M132").Value = -2012.9546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H132").Value = 7853.077
$ws.Range("I132").Value = 14661.6
$ws.Range("J132").Value = 3597.75
$ws.Range("K132").Value = 43984.8
$ws.Range("L132").Value = 10793.25
$ws.Range("M132").Value = -41454.8
$ws.Range("N132").Value = -15853.25

$ws.Range("H138").Value = 34857.25
$ws.Range("J138").Value = 34857.25
$ws.Range("L138").Value = 34857.25
$ws.Range("N138").Value = -45137.25
